# Mark newly implemented opcodes ("i" = Implemented, needs unit tests) on the
# RustyGB opcode-progress sheet. These cells were previously empty (blank /
# "not yet touched") and are now flagged "i" just like their row/col peers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellsToMark = @(
    "D3","H3","L3","P3",
    "C4","D4","H4","L4","P4",
    "C5","D5","F5","G5","H5","L5","P5",
    "H6","P6",
    "H7","P7",
    "H8","P8",
    "B9","C9","D9","E9","F9","G9","I9","P9",
    "H10","P10",
    "H11","P11",
    "H12","P12",
    "H13","P13",
    "H14","P14",
    "H15","P15",
    "H16","P16",
    "H17","P17"
)

foreach ($ref in $cellsToMark) {
    $ws.Range($ref).Value = "i"
}

# Move the active selection to reflect where the author was last working.
$ws.Range("N14").Select()
